# Update contestant status and seat assignments for game show management
#
# 1. Re-order the 3 contestant rows (Felicity, Peter, Kathleen) and flip their
#    Status from "available" to "assigned".
# 2. Insert a new "Seat Assignments" sheet between "Contestants" and "Groups"
#    with the seat-booking rows for the one record day.

$wb = $excel.ActiveWorkbook
$contestants = $wb.Worksheets.Item("Contestants")

# ---------------------------------------------------------------------------
# Step 1: snapshot the existing 3 contestant data rows (row 2..4, cols A..M)
# so we can rewrite them in the new order without losing any values.
# ---------------------------------------------------------------------------
$lastCol = 13  # column M
$snapshot = @{}
for ($r = 2; $r -le 4; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals[$c] = $contestants.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Map old row -> new row (by matching on column A, the ID, so this is not
# order-dependent on how the sheet happened to be laid out):
#   old row 3 (Felicity)  -> new row 2
#   old row 4 (Peter)     -> new row 3
#   old row 2 (Kathleen)  -> new row 4
$byId = @{}
foreach ($r in $snapshot.Keys) {
    $byId[[string]$snapshot[$r][1]] = $snapshot[$r]
}

$newOrder = @(
    "28603f95-d5f6-47ab-88c4-0d79742a6b02",  # Felicity Parker-Hill -> row 2
    "0ccaf8bc-6ade-4ad6-9537-92f9829a6502",  # Peter Adamidis       -> row 3
    "d698b1de-6641-45c6-aa63-f577d2b634bb"   # Kathleen Reynolds    -> row 4
)

$destRow = 2
foreach ($id in $newOrder) {
    $src = $byId[$id]
    for ($c = 1; $c -le $lastCol; $c++) {
        $val = $src[$c]
        if ($c -eq 9) {
            # Column I = Status: available -> assigned
            $val = "assigned"
        }
        $cell = $contestants.Cells.Item($destRow, $c)
        if ($val -eq $null) {
            $cell.ClearContents()
        } else {
            if ($c -eq 6) {
                # Column F = Phone: numeric-looking value that must stay
                # stored as text (matches the source data, which is text).
                $cell.NumberFormat = "@"
            }
            $cell.Value = $val
        }
    }
    $destRow++
}

# ---------------------------------------------------------------------------
# Step 2: insert the new "Seat Assignments" sheet right after "Contestants"
# (i.e. before "Groups"), and populate it.
# ---------------------------------------------------------------------------
$seatSheet = $wb.Worksheets.Add($null, $contestants)
$seatSheet.Name = "Seat Assignments"

$headers = @("ID", "RecordDayID", "ContestantID", "Block", "Seat", "BookingEmailSent", "ConfirmedRSVP", "Notes")
for ($c = 1; $c -le $headers.Length; $c++) {
    $seatSheet.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$recordDayId = "e432f0fe-1383-44a2-990c-5f787da5008a"

$seatRows = @(
    @{ ID = "a001b03c-7c5b-46ba-957e-bb1aeca919d8"; ContestantID = "d698b1de-6641-45c6-aa63-f577d2b634bb"; Block = 1; Seat = "A1" },
    @{ ID = "d5957e66-cb7f-4267-bbb8-1c8c7298b62b"; ContestantID = "28603f95-d5f6-47ab-88c4-0d79742a6b02"; Block = 1; Seat = "A2" },
    @{ ID = "36cd33af-9bac-49bf-bb65-5b0cc9ad1077"; ContestantID = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"; Block = 1; Seat = "A3" }
)

$r = 2
foreach ($row in $seatRows) {
    $seatSheet.Cells.Item($r, 1).Value = $row.ID
    $seatSheet.Cells.Item($r, 2).Value = $recordDayId
    $seatSheet.Cells.Item($r, 3).Value = $row.ContestantID
    $seatSheet.Cells.Item($r, 4).Value = $row.Block
    $seatSheet.Cells.Item($r, 5).Value = $row.Seat
    # BookingEmailSent / ConfirmedRSVP / Notes are left blank (null) for now.
    $r++
}
